$d = $word.ActiveDocument

# Update the date heading (unique text outside the table)
$d.Content.Find.Execute("2024-08-04 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-08-05 Monday", 2)

# Update the division problems in the table, cell by cell, to avoid any
# ambiguity between old/new values that collide with each other.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "54÷3=18, 0" },
    @{ Row = 1;  Col = 2; New = "57÷9=6, 3" },
    @{ Row = 1;  Col = 3; New = "64÷8=8, 0" },
    @{ Row = 1;  Col = 4; New = "52÷7=7, 3" },
    @{ Row = 1;  Col = 5; New = "95÷8=11, 7" },

    @{ Row = 5;  Col = 1; New = "66÷8=8, 2" },
    @{ Row = 5;  Col = 2; New = "78÷7=11, 1" },
    @{ Row = 5;  Col = 3; New = "24÷6=4, 0" },
    @{ Row = 5;  Col = 4; New = "87÷6=14, 3" },
    @{ Row = 5;  Col = 5; New = "82÷8=10, 2" },

    @{ Row = 9;  Col = 1; New = "61÷8=7, 5" },
    @{ Row = 9;  Col = 2; New = "18÷9=2, 0" },
    @{ Row = 9;  Col = 3; New = "18÷3=6, 0" },
    @{ Row = 9;  Col = 4; New = "22÷6=3, 4" },
    @{ Row = 9;  Col = 5; New = "34÷3=11, 1" },

    @{ Row = 13; Col = 1; New = "51÷6=8, 3" },
    @{ Row = 13; Col = 2; New = "86÷4=21, 2" },
    @{ Row = 13; Col = 3; New = "96÷2=48, 0" },
    @{ Row = 13; Col = 4; New = "18÷2=9, 0" },
    @{ Row = 13; Col = 5; New = "10÷8=1, 2" },

    @{ Row = 17; Col = 1; New = "21÷3=7, 0" },
    @{ Row = 17; Col = 2; New = "60÷3=20, 0" },
    @{ Row = 17; Col = 3; New = "34÷9=3, 7" },
    @{ Row = 17; Col = 4; New = "92÷6=15, 2" },
    @{ Row = 17; Col = 5; New = "38÷8=4, 6" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
